# feat: add 2022-Q1 data
#
# Before:  sheets = [ "2021-Q4", "总计" ]
# After:   sheets = [ "2021-Q4", "2022-Q1", "总计" ]
#   - "2022-Q1" is a brand-new sheet (same column layout as "2021-Q4")
#     inserted right after "2021-Q4".
#   - "总计" gets a new row (2022-Q1 totals) inserted above the existing
#     2021-Q4 totals row.

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Add($null, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

# Re-fetch "总计" by name now that the sheet collection has shifted
# (worksheet references resolve by live index, so the handle grabbed
# before the insert would now point at the wrong sheet).
$totalSheet = $wb.Worksheets.Item("总计")

# Copy the header row formatting (bold, centered, bordered) from 2021-Q4
$q4Sheet.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)

$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Copy the index-column formatting (column A header style) down to rows 2-5
$q4Sheet.Range("A2").Copy()
$q1Sheet.Range("A2:A5").PasteSpecial(-4122)

$q1Sheet.Range("A2").Value = 0
$q1Sheet.Range("B2").Value = "'010654"
$q1Sheet.Range("C2").Value = "天弘医药创新混合A"
$q1Sheet.Range("D2").Value = "'11.86"
$q1Sheet.Range("E2").Value = "'87.22"
$q1Sheet.Range("F2").Value = "'3.81"
$q1Sheet.Range("G2").Value = "'0.4519"
$q1Sheet.Range("H2").Value = 9

$q1Sheet.Range("A3").Value = 1
$q1Sheet.Range("B3").Value = "'010054"
$q1Sheet.Range("C3").Value = "万家健康产业混合A"
$q1Sheet.Range("D3").Value = "'8.13"
$q1Sheet.Range("E3").Value = "'86.63"
$q1Sheet.Range("F3").Value = "'3.78"
$q1Sheet.Range("G3").Value = "'0.3073"
$q1Sheet.Range("H3").Value = 5

$q1Sheet.Range("A4").Value = 2
$q1Sheet.Range("B4").Value = "'010655"
$q1Sheet.Range("C4").Value = "天弘医药创新混合C"
$q1Sheet.Range("D4").Value = "'3.52"
$q1Sheet.Range("E4").Value = "'87.22"
$q1Sheet.Range("F4").Value = "'3.81"
$q1Sheet.Range("G4").Value = "'0.1341"
$q1Sheet.Range("H4").Value = 9

$q1Sheet.Range("A5").Value = 3
$q1Sheet.Range("B5").Value = "'010055"
$q1Sheet.Range("C5").Value = "万家健康产业混合C"
$q1Sheet.Range("D5").Value = "'3.36"
$q1Sheet.Range("E5").Value = "'86.63"
$q1Sheet.Range("F5").Value = "'3.78"
$q1Sheet.Range("G5").Value = "'0.1270"
$q1Sheet.Range("H5").Value = 5

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row for 2022-Q1 above the
#    existing 2021-Q4 totals row (which shifts from row 2 to row 3).
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("B2").Value = "'2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 1.02
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1

Write-Output "done"
